$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "98.194.90"
$ws.Range("E2").Value = "  +3.89%  "
$ws.Range("D3").Value = "3.356.25"
$ws.Range("E3").Value = "  +9.12%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "254.87"
$ws.Range("E5").Value = "  +8.53%  "
$ws.Range("D6").Value = "622.24"
$ws.Range("E6").Value = "  +2.19%  "
$ws.Range("E7").Value = "  +8.48%  "
$ws.Range("D8").Value = "0.385"
$ws.Range("E8").Value = "  +2.03%  "
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "3.347.68"
$ws.Range("E10").Value = "  +9.00%  "
$ws.Range("D11").Value = "0.788"
$ws.Range("E11").Value = "  -1.92%  "
$ws.Range("E12").Value = "  +1.61%  "
$ws.Range("D13").Value = "98.012.96"
$ws.Range("E13").Value = "  +4.18%  "
$ws.Range("D14").Value = "35.86"
$ws.Range("E14").Value = "  +6.43%  "
$ws.Range("D15").Value = "0.0000245"
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("E16").Value = "  +9.30%  "
$ws.Range("D17").Value = "5.48"
$ws.Range("E17").Value = "  +2.98%  "
$ws.Range("D18").Value = "3.364.10"
$ws.Range("E18").Value = "  +10.63%  "
$ws.Range("D19").Value = "3.58"
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("D20").Value = "14.86"
$ws.Range("E20").Value = "  +3.27%  "
$ws.Range("D21").Value = "483.60"
$ws.Range("E21").Value = "  +10.02%  "
$ws.Range("D22").Value = "5.83"
$ws.Range("E22").Value = "  +3.06%  "
$ws.Range("D23").Value = "0.0000206"
$ws.Range("E23").Value = "  +9.58%  "
$ws.Range("D24").Value = "9.07"
$ws.Range("E24").Value = "  +3.02%  "
$ws.Range("D25").Value = "5.65"
$ws.Range("E25").Value = "  +2.86%  "
$ws.Range("D26").Value = "87.75"
$ws.Range("E26").Value = "  +3.62%  "
$ws.Range("D27").Value = "11.95"
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("D28").Value = "3.540.43"
$ws.Range("E28").Value = "  +9.71%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  +5.09%  "
$ws.Range("D31").Value = "0.243"
$ws.Range("E31").Value = "  -1.17%  "
$ws.Range("E32").Value = "  -2.00%  "
$ws.Range("E33").Value = "  +12.86%  "
$ws.Range("D34").Value = "9.19"
$ws.Range("E34").Value = "  +1.79%  "
$ws.Range("D35").Value = "27.14"
$ws.Range("E35").Value = "  +6.93%  "
$ws.Range("D36").Value = "0.151"
$ws.Range("E36").Value = "  -1.90%  "
$ws.Range("D37").Value = "511.43"
$ws.Range("E37").Value = "  +9.25%  "
$ws.Range("D38").Value = "7.27"
$ws.Range("E38").Value = "  -4.59%  "
$ws.Range("D39").Value = "1.94"
$ws.Range("E39").Value = "  +4.17%  "
$ws.Range("D40").Value = "24.85"
$ws.Range("E40").Value = "  +3.48%  "
$ws.Range("D41").Value = "0.445"
$ws.Range("E41").Value = "  +1.93%  "
$ws.Range("D42").Value = "1.26"
$ws.Range("E42").Value = "  +0.55%  "
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("D44").Value = "3.22"
$ws.Range("E44").Value = "  +4.24%  "
$ws.Range("D45").Value = "0.782"
$ws.Range("E45").Value = "  +16.89%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "159.66"
$ws.Range("E47").Value = "  -1.18%  "
$ws.Range("D48").Value = "1.91"
$ws.Range("E48").Value = "  +5.43%  "
$ws.Range("D49").Value = "1.36"
$ws.Range("E49").Value = "  +6.82%  "
$ws.Range("D50").Value = "45.46"
$ws.Range("E50").Value = "  +4.40%  "
$ws.Range("D51").Value = "4.49"
$ws.Range("E51").Value = "  +5.89%  "

$ws.Range("D2:D51").Style = "Normal"
